# Tidsregistrering i PTE projektet Simon Nielsen.xlsx
# Har lavet Test for Boejningsspaending og opdateret tidsregistrering
#
# Adds 5 new time-registration rows (52-56) for 2017-04-22 on the
# "Tidsregistrering" sheet, reusing the existing date/time number formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

function Add-Entry([int]$row, [object]$date, [string]$text, [double]$start, [double]$end) {
    if ($date -ne $null) {
        $ws.Range("A41").Copy() | Out-Null
        $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
        $ws.Range("A$row").Value = $date
    }

    $ws.Range("F$row").Value = $text

    $ws.Range("G41").Copy() | Out-Null
    $ws.Range("G$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("G$row").Value = $start

    $ws.Range("H41").Copy() | Out-Null
    $ws.Range("H$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("H$row").Value = $end
}

Add-Entry 52 42847 "Rettet Test Suite OC15"                      0.34722222222222227 0.36458333333333331
Add-Entry 53 $null "Implementeret OC12"                          0.34722222222222227 0.38541666666666669
Add-Entry 54 $null "Kode Review OC13"                            0.3923611111111111  0.4201388888888889
Add-Entry 55 $null "Hjalp med implementering af OC 15 & 16"      0.42708333333333331 0.47916666666666669
Add-Entry 56 $null "Test test og test"                           0.51041666666666663 0.64583333333333337

$ws.Range("I57").Select() | Out-Null
